# Adding files and making updates to framework for both POST and GET methods
#
# Populates row 3 of the "TestData" sheet as a second scenario ("Scenario 2",
# a GET endpoint) mirroring row 2's layout, wires up a hyperlink on E3 like
# the one already on E2, and moves the active selection down to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- 1) Row 3 values: same shape as row 2, new scenario name + GET endpoint ---
# (new shared strings are appended in the order they are first written, so
# write D3 before A3 to match the endpoint-then-scenario ordering used by
# the target workbook.)
$ws.Range("D3").Value2 = "/maps/api/place/get/json"
$ws.Range("A3").Value2 = "Scenario 2"
$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2
$ws.Range("H3").Value2 = -38.383493999999999
$ws.Range("I3").Value2 = 33.427362000000002
$ws.Range("J3").Value2 = 50
$ws.Range("K3").Value2 = $ws.Range("K2").Value2
$ws.Range("L3").Value2 = $ws.Range("L2").Value2
$ws.Range("M3").Value2 = $ws.Range("M2").Value2
$ws.Range("N3").Value2 = $ws.Range("N2").Value2
$ws.Range("O3").Value2 = $ws.Range("O2").Value2
$ws.Range("P3").Value2 = $ws.Range("P2").Value2

# --- 2) Hyperlink on E3, mirroring the existing one on E2 ---
$ws.Hyperlinks.Add($ws.Range("E3"), "https://rahulshettyacademy.com/")

# --- 3) Re-apply row 2's cell formatting (Value2/hyperlink writes above ---
#        reset formatting on the touched cells, so copy row 2's look back
#        onto row 3 now that the content is in place).
$ws.Range("A2:P2").Copy()
$ws.Range("A3:P3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- 4) Move the view/selection down to the new row ---
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("A3").Select()
